$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 47 (shifts existing rows 47-58 down to 48-59)
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new menu item
$ws.Range("A47").Value = "Snacks"
$ws.Range("B47").Value = "Melody toffee 1 rs"
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = "Melody toffee 1 rs.jpg"
$ws.Range("F47").Value = "Fast Food"

# Update the active selection to match the authored change
$ws.Range("F47").Select()

# The filter database defined name grows by one row to include the new row
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$F`$57"

$wb.Save()
